# Import PG baru, Up Author
# Insert a new "Kelompok Passing Grade" column between the existing
# "Prodi" (A) and "Passing Grade" (B) columns, shifting the latter to C.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Passing Grade" column from B to C and create the
# new column B for "Kelompok Passing Grade".
$ws.Columns("B").Insert()

$ws.Range("B1").Value = "Kelompok Passing Grade"

# Match the column's width as closely as this runtime's width grid allows.
$ws.Columns("B").ColumnWidth = 24.67

# Reflect the author's final cell selection.
$ws.Range("C8").Select()
